{"js": "// Insert \"(4)\" (the number of people who registered) into the two\n// \"Low number of sign up's\" follow-up bullets (English + Swedish) under\n// \"\u00c5re Skiweek 2007\", matching:\n//   \"Low number of sign up's...\"        -> \"Low number (4) of sign up's...\"\n//   \"Har haft l\u00e5gt antal anm\u00e4lda...\"    -> \"Har haft l\u00e5gt antal (4) anm\u00e4lda...\"\n//\n// The source text is split into three runs (prefix / \"(4)\" / suffix) that\n// all keep the original bold run formatting, mirroring how Word splits a\n// run when you type new text in the middle of it.\n\nconst body = context.document.body;\n\nasync function insertNumberMarker(anchorText, fullSentenceNeedle) {\n  // Locate the unique anchor (\"...word \" right before the part that gets\n  // split) and insert the marker + trailing space right after it.\n  const anchorResults = body.search(anchorText, { matchCase: true, matchWholeWord: false });\n  await context.sync();\n\n  if (anchorResults.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for anchor ${JSON.stringify(anchorText)}, got ${anchorResults.items.length}`\n    );\n  }\n\n  const anchorRange = anchorResults.items[0];\n  anchorRange.insertText(\"(4) \", Word.InsertLocation.after);\n  await context.sync();\n\n  // Re-search, scoped to the sentence we just edited, so we grab the right\n  // \"(4)\" occurrence (English and Swedish sentences both get one).\n  const sentenceResults = body.search(fullSentenceNeedle, { matchCase: true, matchWholeWord: false });\n  await context.sync();\n  if (sentenceResults.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for sentence fragment ${JSON.stringify(fullSentenceNeedle)}, got ${sentenceResults.items.length}`\n    );\n  }\n  const sentenceRange = sentenceResults.items[0];\n\n  const markerResults = sentenceRange.search(\"(4)\", { matchCase: true, matchWholeWord: false });\n  markerResults.load(\"font\");\n  await context.sync();\n  if (markerResults.items.length !== 1) {\n    throw new Error(`Expected exactly one \"(4)\" inside the sentence, got ${markerResults.items.length}`);\n  }\n  const markerRange = markerResults.items[0];\n\n  // Toggling bold off/on forces the run containing \"(4)\" to become its own\n  // run (splitting it away from the prefix/suffix text), while leaving the\n  // final bold value identical to the surrounding text (no stray overrides).\n  markerRange.font.bold = false;\n  await context.sync();\n  markerRange.font.bold = true;\n  await context.sync();\n}\n\n// English bullet: \"Low number of sign up's. Will fix a promotion picture...\"\nawait insertNumberMarker(\"Low number \", \"Low number (4) of sign up\");\n\n// Swedish bullet: \"Har haft l\u00e5gt antal anm\u00e4lda. Ska fixa bild...\"\nawait insertNumberMarker(\"Har haft l\u00e5gt antal \", \"Har haft l\u00e5gt antal (4) anm\u00e4lda\");\n", "ps1": "# Insert \"(4)\" (the number of people who registered) into the two\n# \"Low number of sign up's\" follow-up bullets (English + Swedish) under\n# \"\u00c5re Skiweek 2007\", matching:\n#   \"Low number of sign up's...\"        -> \"Low number (4) of sign up's...\"\n#   \"Har haft l\u00e5gt antal anm\u00e4lda...\"    -> \"Har haft l\u00e5gt antal (4) anm\u00e4lda...\"\n#\n# The source text is split into three runs (prefix / \"(4)\" / suffix) that\n# all keep the original bold run formatting, mirroring how Word splits a\n# run when you type new text in the middle of it.\n\n$d = $word.ActiveDocument\n\nfunction Insert-NumberMarker($anchorText) {\n    # Locate the unique anchor (\"...word \" right before the part that gets\n    # split), collapse to its end, and insert the marker + trailing space.\n    $searchRange = $d.Content\n    $find = $searchRange.Find\n    $find.Text = $anchorText\n    $find.MatchCase = $true\n    $find.Execute() | Out-Null\n    if (-not $find.Found) {\n        throw \"Anchor text not found: $anchorText\"\n    }\n\n    $searchRange.Collapse(0)  # wdCollapseEnd\n    $searchRange.InsertAfter(\"(4) \")\n\n    # Find \"(4)\" again, scoped to start right at the insertion point so we\n    # can't accidentally match an earlier \"(4)\" already written elsewhere\n    # (the English bullet is processed before the Swedish one).\n    $markerRange = $d.Range($searchRange.Start, $d.Content.End)\n    $markerFind = $markerRange.Find\n    $markerFind.Text = \"(4)\"\n    $markerFind.MatchCase = $true\n    $markerFind.Execute() | Out-Null\n    if (-not $markerFind.Found) {\n        throw \"Inserted (4) marker not found after anchor: $anchorText\"\n    }\n\n    # Toggling bold off/on forces the run containing \"(4)\" to become its own\n    # run (splitting it away from the prefix/suffix text), while leaving the\n    # final bold value identical to the surrounding text (no stray overrides).\n    $markerRange.Font.Bold = 0\n    $markerRange.Font.Bold = 1\n}\n\n# English bullet: \"Low number of sign up's. Will fix a promotion picture...\"\nInsert-NumberMarker \"Low number \"\n\n# Swedish bullet: \"Har haft l\u00e5gt antal anm\u00e4lda. Ska fixa bild...\"\nInsert-NumberMarker \"Har haft l\u00e5gt antal \"\n"}
